# Sync attendance_reports: reorder "Recorded By" (column G) entries so that
# the real/flagged account email (dnasr281@gmail.com or backup@backdoor.com)
# is listed before the generic "System" entries, instead of after them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Emails that should be moved to the front of the "Recorded By" list when
# they appear alongside a leading "System" entry.
$priorityEmails = @('dnasr281@gmail.com', 'backup@backdoor.com')

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value()

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -eq '') { continue }

    $parts = $text -split ',\s*'

    if ($parts.Count -lt 2) { continue }
    if ($parts[0] -ne 'System') { continue }

    # Only reorder when one of the remaining parts is a priority email and
    # none of the non-"System" parts fall outside the priority list (e.g.
    # leave "System, admin@admin.com" untouched).
    $rest = $parts[1..($parts.Count - 1)]
    $nonSystemRest = $rest | Where-Object { $_ -ne 'System' -and $_ -ne 'system' }

    if ($nonSystemRest.Count -eq 0) { continue }

    $matchesPriority = $true
    foreach ($item in $nonSystemRest) {
        if ($priorityEmails -notcontains $item) {
            $matchesPriority = $false
            break
        }
    }
    if (-not $matchesPriority) { continue }

    # Move the last element to the front, preserving the relative order of
    # the remaining entries (e.g. "System, system, backup@backdoor.com"
    # becomes "backup@backdoor.com, System, system").
    $newOrder = @($parts[-1]) + $parts[0..($parts.Count - 2)]
    $newText = [string]::Join(', ', $newOrder)

    $cell.Value = $newText
}
